# Append one new data row (row 85) to the stream-stats table, continuing
# the existing stream/avg/max/follow series (8-5 data point), matching
# the formatting of the preceding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 84
$newRow = 85

# Copy formatting (style) from the last existing data row onto the new row.
$ws.Range("A$lastRow`:D$lastRow").Copy()
$ws.Range("A$newRow`:D$newRow").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new data values.
$ws.Cells.Item($newRow, 1).Value = 83.0
$ws.Cells.Item($newRow, 2).Value = 179.0
$ws.Cells.Item($newRow, 3).Value = 214.0
$ws.Cells.Item($newRow, 4).Value = 67.0
